$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell "time_taken" in F1, matching the style of the
# existing header row (bold font, border, centered) by copying E1's format.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Populate F2:F62 with the recorded time_taken timestamps (plain text,
# no special style - matches the rest of the data rows).
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:51:30.912042"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:51:30.912055"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:51:30.912059"
$ws.Cells.Item(5, 6).Value = "2021-10-05 10:51:30.912062"
$ws.Cells.Item(6, 6).Value = "2021-10-05 10:51:30.912066"
$ws.Cells.Item(7, 6).Value = "2021-10-05 10:51:30.912069"
$ws.Cells.Item(8, 6).Value = "2021-10-05 10:51:30.912072"
$ws.Cells.Item(9, 6).Value = "2021-10-05 10:51:30.912075"
$ws.Cells.Item(10, 6).Value = "2021-10-05 10:51:30.912078"
$ws.Cells.Item(11, 6).Value = "2021-10-05 10:51:30.912081"
$ws.Cells.Item(12, 6).Value = "2021-10-05 10:51:30.912084"
$ws.Cells.Item(13, 6).Value = "2021-10-05 10:51:30.912087"
$ws.Cells.Item(14, 6).Value = "2021-10-05 10:51:30.912090"
$ws.Cells.Item(15, 6).Value = "2021-10-05 10:51:30.912093"
$ws.Cells.Item(16, 6).Value = "2021-10-05 10:51:30.912096"
$ws.Cells.Item(17, 6).Value = "2021-10-05 10:51:30.912099"
$ws.Cells.Item(18, 6).Value = "2021-10-05 10:51:30.912102"
$ws.Cells.Item(19, 6).Value = "2021-10-05 10:51:30.912106"
$ws.Cells.Item(20, 6).Value = "2021-10-05 10:51:30.912109"
$ws.Cells.Item(21, 6).Value = "2021-10-05 10:51:30.912112"
$ws.Cells.Item(22, 6).Value = "2021-10-05 10:51:30.912115"
$ws.Cells.Item(23, 6).Value = "2021-10-05 10:51:30.912118"
$ws.Cells.Item(24, 6).Value = "2021-10-05 10:51:30.912120"
$ws.Cells.Item(25, 6).Value = "2021-10-05 10:51:30.912123"
$ws.Cells.Item(26, 6).Value = "2021-10-05 10:51:30.912127"
$ws.Cells.Item(27, 6).Value = "2021-10-05 10:51:30.912130"
$ws.Cells.Item(28, 6).Value = "2021-10-05 10:51:30.912133"
$ws.Cells.Item(29, 6).Value = "2021-10-05 10:51:30.912136"
$ws.Cells.Item(30, 6).Value = "2021-10-05 10:51:30.912139"
$ws.Cells.Item(31, 6).Value = "2021-10-05 10:51:30.912142"
$ws.Cells.Item(32, 6).Value = "2021-10-05 10:51:30.912145"
$ws.Cells.Item(33, 6).Value = "2021-10-05 10:51:30.912148"
$ws.Cells.Item(34, 6).Value = "2021-10-05 10:51:30.912151"
$ws.Cells.Item(35, 6).Value = "2021-10-05 10:51:30.912154"
$ws.Cells.Item(36, 6).Value = "2021-10-05 10:51:30.912157"
$ws.Cells.Item(37, 6).Value = "2021-10-05 10:51:30.912160"
$ws.Cells.Item(38, 6).Value = "2021-10-05 10:51:30.912163"
$ws.Cells.Item(39, 6).Value = "2021-10-05 10:51:30.912166"
$ws.Cells.Item(40, 6).Value = "2021-10-05 10:51:30.912169"
$ws.Cells.Item(41, 6).Value = "2021-10-05 10:51:30.912172"
$ws.Cells.Item(42, 6).Value = "2021-10-05 10:51:30.912176"
$ws.Cells.Item(43, 6).Value = "2021-10-05 10:51:30.912179"
$ws.Cells.Item(44, 6).Value = "2021-10-05 10:51:30.912182"
$ws.Cells.Item(45, 6).Value = "2021-10-05 10:51:30.912185"
$ws.Cells.Item(46, 6).Value = "2021-10-05 10:51:30.912188"
$ws.Cells.Item(47, 6).Value = "2021-10-05 10:51:30.912191"
$ws.Cells.Item(48, 6).Value = "2021-10-05 10:51:30.912194"
$ws.Cells.Item(49, 6).Value = "2021-10-05 10:51:30.912197"
$ws.Cells.Item(50, 6).Value = "2021-10-05 10:51:30.912200"
$ws.Cells.Item(51, 6).Value = "2021-10-05 10:51:30.912203"
$ws.Cells.Item(52, 6).Value = "2021-10-05 10:51:30.912206"
$ws.Cells.Item(53, 6).Value = "2021-10-05 10:51:30.912209"
$ws.Cells.Item(54, 6).Value = "2021-10-05 10:51:30.912212"
$ws.Cells.Item(55, 6).Value = "2021-10-05 10:51:30.912215"
$ws.Cells.Item(56, 6).Value = "2021-10-05 10:51:30.912218"
$ws.Cells.Item(57, 6).Value = "2021-10-05 10:51:30.912221"
$ws.Cells.Item(58, 6).Value = "2021-10-05 10:51:30.912224"
$ws.Cells.Item(59, 6).Value = "2021-10-05 10:51:30.912227"
$ws.Cells.Item(60, 6).Value = "2021-10-05 10:51:30.912230"
$ws.Cells.Item(61, 6).Value = "2021-10-05 10:51:30.912233"
$ws.Cells.Item(62, 6).Value = "2021-10-05 10:51:30.912236"